$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 379.3
$ws.Range("I11").Value = 379.3
$ws.Range("K11").Value = 379.3
$ws.Range("M11").Value = -239.3
$ws.Range("H15").Value = 806.6818
$ws.Range("I15").Value = 806.6818
$ws.Range("K15").Value = 2420.0454
$ws.Range("M15").Value = -2251.0454
$ws.Range("H17").Value = 441343.78
$ws.Range("J17").Value = 472842.03
$ws.Range("L17").Value = 1418526.09
$ws.Range("N17").Value = -1418862.09
$ws.Range("H33").Value = 123668.78
$ws.Range("I33").Value = 139077.38
$ws.Range("K33").Value = 139077.38
$ws.Range("M33").Value = -138848.38
$ws.Range("H86").Value = 5649.6
$ws.Range("I86").Value = 7156.4287
$ws.Range("J86").Value = 2133.6667
$ws.Range("K86").Value = 7156.4287
$ws.Range("L86").Value = 2133.6667
$ws.Range("M86").Value = -6033.4287
$ws.Range("N86").Value = -4379.6667
$ws.Range("H89").Value = 5649.6
$ws.Range("I89").Value = 7156.4287
$ws.Range("J89").Value = 2133.6667
$ws.Range("K89").Value = 35782.14350000001
$ws.Range("L89").Value = 10668.3335
$ws.Range("M89").Value = -30166.14350000001
$ws.Range("N89").Value = -21900.3335
$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680
$ws.Range("H113").Value = 3353.111
$ws.Range("I113").Value = 3353.111
$ws.Range("K113").Value = 3353.111
$ws.Range("M113").Value = -99.11099999999988
$ws.Range("H116").Value = 7625.7915
$ws.Range("I116").Value = 7625.7915
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 7625.7915
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -4183.7915
$ws.Range("N116").Value = ""
$ws.Range("H132").Value = 2488.3333
$ws.Range("I132").Value = 2532.7273
$ws.Range("K132").Value = 7598.1819
$ws.Range("M132").Value = -5068.1819
$ws.Range("H138").Value = 3004.6924
$ws.Range("I138").Value = 2088.4666
$ws.Range("J138").Value = 4254.091
$ws.Range("K138").Value = 6265.399800000001
$ws.Range("L138").Value = 12762.273
$ws.Range("M138").Value = -1125.399800000001
$ws.Range("N138").Value = -23042.273

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4241.231
$ws.Range("I32").Value = 3055.7778
$ws.Range("K32").Value = 3055.7778
$ws.Range("M32").Value = -2768.7778
$ws.Range("H45").Value = 6119.8823
$ws.Range("I45").Value = 6817.6924
$ws.Range("K45").Value = 6817.6924
$ws.Range("M45").Value = -6440.6924
$ws.Range("H61").Value = 43480536
$ws.Range("I61").Value = 66668050
$ws.Range("J61").Value = 3956.75
$ws.Range("K61").Value = 66668050
$ws.Range("L61").Value = 3956.75
$ws.Range("M61").Value = -66667838
$ws.Range("N61").Value = -4380.75
$ws.Range("H110").Value = 334732.66
$ws.Range("I110").Value = 501199.5
$ws.Range("K110").Value = 501199.5
$ws.Range("M110").Value = -499154.5
$ws.Range("H122").Value = 7489.9165
$ws.Range("I122").Value = 6208.778
$ws.Range("K122").Value = 18626.334
$ws.Range("M122").Value = -16176.334
$ws.Range("H132").Value = 7145938
$ws.Range("I132").Value = 8335427.5
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 25006282.5
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -25003752.5
$ws.Range("N132").Value = -32060
$ws.Range("H136").Value = 43480536
$ws.Range("I136").Value = 66668050
$ws.Range("J136").Value = 3956.75
$ws.Range("K136").Value = 200004150
$ws.Range("L136").Value = 11870.25
$ws.Range("M136").Value = -200001600
$ws.Range("N136").Value = -16970.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3059.6
$ws.Range("I105").Value = 2119.2
$ws.Range("K105").Value = 2119.2
$ws.Range("M105").Value = -372.1999999999998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""
$ws.Range("H58").Value = 17246994
$ws.Range("I58").Value = 41676276
$ws.Range("K58").Value = 41676276
$ws.Range("M58").Value = -41676073
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H74").Value = 46651.332
$ws.Range("J74").Value = 49977
$ws.Range("L74").Value = 49977
$ws.Range("N74").Value = -51725
$ws.Range("H77").Value = 46651.332
$ws.Range("J77").Value = 49977
$ws.Range("L77").Value = 149931
$ws.Range("N77").Value = -158667
$ws.Range("H107").Value = 551463.6
$ws.Range("I107").Value = 1010786.8
$ws.Range("K107").Value = 1010786.8
$ws.Range("M107").Value = -1008866.8
$ws.Range("H120").Value = 78999.5
$ws.Range("J120").Value = 78999.5
$ws.Range("L120").Value = 78999.5
$ws.Range("N120").Value = -86257.5
$ws.Range("H132").Value = 250001520
$ws.Range("I132").Value = 250001520
$ws.Range("K132").Value = 750004560
$ws.Range("M132").Value = -750002030
$ws.Range("H136").Value = 17246994
$ws.Range("I136").Value = 41676276
$ws.Range("K136").Value = 125028828
$ws.Range("M136").Value = -125026278

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 12227.217
$ws.Range("I56").Value = 12227.217
$ws.Range("K56").Value = 12227.217
$ws.Range("M56").Value = -11697.217
$ws.Range("H107").Value = 859.6667
$ws.Range("I107").Value = 191.75
$ws.Range("K107").Value = 575.25
$ws.Range("M107").Value = 1344.75
$ws.Range("H131").Value = 2093.68
$ws.Range("J131").Value = 2212.5
$ws.Range("L131").Value = 6637.5
$ws.Range("N131").Value = -16717.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 53499
$ws.Range("J15").Value = 53499
$ws.Range("L15").Value = 53499
$ws.Range("N15").Value = -54075
$ws.Range("H81").Value = 53499
$ws.Range("J81").Value = 53499
$ws.Range("L81").Value = 53499
$ws.Range("N81").Value = -55495
$ws.Range("H84").Value = 53499
$ws.Range("J84").Value = 53499
$ws.Range("L84").Value = 160497
$ws.Range("N84").Value = -170481
$ws.Range("H102").Value = 3451.2354
$ws.Range("I102").Value = 3047.3125
$ws.Range("J102").Value = 9914
$ws.Range("K102").Value = 3047.3125
$ws.Range("L102").Value = 9914
$ws.Range("M102").Value = -1425.3125
$ws.Range("N102").Value = -13158
$ws.Range("H122").Value = 127764.5
$ws.Range("I122").Value = 155455.62
$ws.Range("J122").Value = 17000
$ws.Range("K122").Value = 466366.86
$ws.Range("L122").Value = 51000
$ws.Range("M122").Value = -463916.86
$ws.Range("N122").Value = -55900
$ws.Range("H126").Value = 4160.769
$ws.Range("I126").Value = 2509
$ws.Range("J126").Value = 9666.666999999999
$ws.Range("K126").Value = 7527
$ws.Range("L126").Value = 29000.001
$ws.Range("M126").Value = -5057
$ws.Range("N126").Value = -33940.001
$ws.Range("H132").Value = 8931156
$ws.Range("I132").Value = 9617945
$ws.Range("K132").Value = 28853835
$ws.Range("M132").Value = -28851305
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = ""

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3841.4443
$ws.Range("I40").Value = 3696.75
$ws.Range("K40").Value = 3696.75
$ws.Range("M40").Value = -3560.75
$ws.Range("H46").Value = 2091.5
$ws.Range("I46").Value = 2091.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2091.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1903.5
$ws.Range("N46").Value = ""
$ws.Range("H122").Value = 13450.5
$ws.Range("I122").Value = 12414.167
$ws.Range("K122").Value = 37242.501
$ws.Range("M122").Value = -34792.501
$ws.Range("H132").Value = 53337600
$ws.Range("I132").Value = 60003924
$ws.Range("K132").Value = 180011772
$ws.Range("M132").Value = -180009242
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -17100

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1165.8948
$ws.Range("I113").Value = 1064
$ws.Range("K113").Value = 3192
$ws.Range("M113").Value = -1022
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = ""
$ws.Range("H126").Value = 2466.3333
$ws.Range("I126").Value = 2499.75
$ws.Range("K126").Value = 7499.25
$ws.Range("M126").Value = -5029.25
$ws.Range("H132").Value = 12828685
$ws.Range("I132").Value = 18525300
$ws.Range("K132").Value = 55575900
$ws.Range("M132").Value = -55573370
$ws.Range("H136").Value = 26317776
$ws.Range("I136").Value = 31252034
$ws.Range("K136").Value = 93756102
$ws.Range("M136").Value = -93753552
